$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}

Replace-Text "2024-01-13 Saturday" "2024-01-14 Sunday"

Replace-Text "81×18=1458" "30×75=2250"
Replace-Text "54×69=3726" "84×41=3444"
Replace-Text "30×90=2700" "86×82=7052"
Replace-Text "56×92=5152" "87×95=8265"
Replace-Text "99×74=7326" "48×54=2592"

Replace-Text "82×97=7954" "45×64=2880"
Replace-Text "62×65=4030" "55×97=5335"
Replace-Text "68×18=1224" "76×33=2508"
Replace-Text "98×25=2450" "75×13=975"
Replace-Text "94×80=7520" "52×15=780"

Replace-Text "92×66=6072" "65×84=5460"
Replace-Text "23×62=1426" "33×64=2112"
Replace-Text "74×35=2590" "36×25=900"
Replace-Text "21×12=252" "82×60=4920"
Replace-Text "60×36=2160" "17×91=1547"

Replace-Text "42×68=2856" "15×14=210"
Replace-Text "45×29=1305" "97×55=5335"
Replace-Text "78×53=4134" "28×94=2632"
Replace-Text "65×67=4355" "62×52=3224"
Replace-Text "22×14=308" "36×45=1620"

Replace-Text "92×73=6716" "15×79=1185"
Replace-Text "70×94=6580" "84×96=8064"
Replace-Text "87×42=3654" "69×88=6072"
Replace-Text "68×51=3468" "16×46=736"
Replace-Text "76×38=2888" "75×41=3075"
